$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text (e.g. "28.628.45"),
# even though many of the updated quotes look like ordinary numbers (e.g.
# "317.09"). Force those specific cells to Text format first so Excel
# keeps them as literal strings instead of silently parsing them as numbers.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D14", "D15", "D17", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed cryptocurrency price/volume data (rows 40-42 also
# reorder which coin occupies which rank, per the source diff).
$ws.Range("D2").Value = '28.599.83'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.803.57'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '317.09'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.5335'
$ws.Range("E7").Value = '  -6.66%  '
$ws.Range("D8").Value = '0.3759'
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").Value = '0.07504'
$ws.Range("E9").Value = '  -1.61%  '
$ws.Range("D10").Value = '42.39'
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("E13").Value = '  -2.58%  '
$ws.Range("D14").Value = '6.143'
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("D15").Value = '7.373'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '1.799.53'
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").Value = '90.18'
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '0.06452'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  -0.56%  '
$ws.Range("D22").Value = '5.915'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").Value = '28.637.69'
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").Value = '11.09'
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("D25").Value = '2.098'
$ws.Range("E25").Value = '  -1.23%  '
$ws.Range("D26").Value = '158.46'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '20.43'
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").Value = '2.008.62'
$ws.Range("E28").Value = '  -0.69%  '
$ws.Range("D29").Value = '2.345'
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("D30").Value = '122.79'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = '1.103'
$ws.Range("E31").Value = '  -5.96%  '
$ws.Range("D32").Value = '0.1062'
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").Value = '5.640'
$ws.Range("E33").Value = '  -2.66%  '
$ws.Range("D34").Value = '3.684'
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("D35").Value = '0.2245'
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("D36").Value = '0.06459'
$ws.Range("E36").Value = '  +5.91%  '
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '8.747'
$ws.Range("E38").Value = '  -2.55%  '
$ws.Range("D39").Value = '5.037'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.24'
$ws.Range("E40").Value = '  -4.26%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6216'
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.193'
$ws.Range("E42").Value = '  +3.26%  '
$ws.Range("D43").Value = '1.426'
$ws.Range("E43").Value = '  +3.50%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '13.31'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").Value = '3.689'
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '0.5847'
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("D48").Value = '126.31'
$ws.Range("E48").Value = '  +3.28%  '
$ws.Range("D49").Value = '1.939'
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").Value = '1.153'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '0.06893'
$ws.Range("E51").Value = '  +0.53%  '
